$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for column F (dSF) — repulled data / push all data / mean calculation
$updates = @{
    3  = 8
    4  = -4
    9  = -3
    11 = 3
    13 = -1
    17 = 1
    24 = -1
    26 = 3
    29 = 2
    30 = 1
    32 = 2
    37 = 2
    38 = 3
    40 = 1
    43 = -1
    44 = -1
    45 = 4
    46 = 2
    47 = -1
    54 = -1
    55 = 5
    57 = 0
    60 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
